$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set the Runmode column (C) to "Y" for rows 3-7, matching row 2's value
$ws.Range("C3:C7").Value = "Y"

# Update the selection to reflect the newly run rows
$ws.Range("C2:C7").Select()
